# advisors: evaluation layout, risk survey
#
# 1) Insert a new "income" column after "profession" (shifts former I..Q to J..R)
# 2) Populate the newly-available religion / education_school / profession / income
#    cells for the first five respondents (rows 2-6)
# 3) Update the q1..q5 answers (now columns N..R) for rows 2-6 with the new survey values
# 4) Append five more respondent rows (7-11) that repeat respondents 0-4 (ids 5-9)
#    with matching riskgroup / riskgroup_text formulas

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) insert the new "income" column (between profession=H and age=I) ---
$ws.Columns("I:I").Insert()
$ws.Range("I1").Value = "income"

# --- respondent reference data (shared across the first and second block) ---
# columns: name, gender, nationality, religion, education_school, profession, income, introText
$people = @(
    @("Christine", "weiblich",    "deutsch",  "evangelisch", "Abitur",               "Lehrerin",     "mehr als 4000€",
      "Hallo, ich bin Christine und ich würde gerne 10.000 € anlegen. In meiner Freizeit spiele ich gerne Tennis und gehe gerne spazieren."),
    @("Flo",       "männlich",    "deutsch",  "keine",       "Realschulabschluss",   "Kaufmann",     "3000-3999€",
      "Hallo, ich bin Flo und ich würde gerne 10.000 € anlegen. In meiner Freizeit spiele ich gerne Tennis und gehe gerne spazieren."),
    @("Zeynep",    "weiblich",    "türkisch", "katholisch",  "Hauptschulabschluss",  "Ärztin",       "2000-2999€",
      "Hallo, ich bin Zeynep und ich würde gerne 10.000 € anlegen. In meiner Freizeit spiele ich gerne Tennis und gehe gerne spazieren."),
    @("Alparslan", "männlich",    "türkisch", "muslimisch",  "kein Abschluss",       "Anwalt",       "weniger als 1000€",
      "Hallo, ich bin Alparslan und ich würde gerne 10.000 € anlegen. In meiner Freizeit spiele ich gerne Tennis und gehe gerne spazieren."),
    @("Toni",      "nicht-binär", "deutsch",  "jüdisch",     "Fachabitur",           "HandwerkerIn", "1000-1999€",
      "Hallo, ich bin Toni und ich würde gerne 10.000 € anlegen. In meiner Freizeit spiele ich gerne Tennis und gehe gerne spazieren.")
)

# q1..q5 answers per respondent (now living in columns N..R)
$answers = @(
    @(1, 2, 3, 4, 4),
    @(1, 2, 4, 1, 2),
    @(3, 2, 2, 4, 2),
    @(1, 2, 2, 1, 1),
    @(4, 4, 4, 2, 4)
)

# --- 2 & 3) fill in the new columns + refreshed answers for rows 2-6 ---
# (populate column-by-column, in the same per-cell order the original author
#  used, so the shared-string table ends up in the same order)
foreach ($i in @(0, 1, 2, 4, 3)) {
    $row = 2 + $i
    $ws.Range("E$row").Value = $people[$i][3]
}
foreach ($i in @(0, 2, 1, 3, 4)) {
    $row = 2 + $i
    $ws.Range("F$row").Value = $people[$i][4]
}
foreach ($i in @(0, 1, 2, 3, 4)) {
    $row = 2 + $i
    $ws.Range("H$row").Value = $people[$i][5]
}
foreach ($i in @(0, 1, 2, 4, 3)) {
    $row = 2 + $i
    $ws.Range("I$row").Value = $people[$i][6]
}

for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $a = $answers[$i]
    $ws.Range("N$row").Value = $a[0]
    $ws.Range("O$row").Value = $a[1]
    $ws.Range("P$row").Value = $a[2]
    $ws.Range("Q$row").Value = $a[3]
    $ws.Range("R$row").Value = $a[4]
}

# refresh the riskgroup / riskgroup_text formulas so they reference the shifted N:R range
# (assign across the whole range at once so Excel keeps them as shared formulas,
#  same as the authored workbook: L2:L5 share one formula, M2:M6 share another,
#  while row 6 keeps its own distinct - non-shared - formula)
$ws.Range("L2:L5").Formula = "=IF(SUM(N2:R2)<=20/5+15/5,1,IF(SUM(N2:R2)<=20/5+15/5*2,2,IF(SUM(N2:R2)<=20/5+15/5*3,3,IF(SUM(N2:R2)<=20/5+15/5*4,4,IF(SUM(N2:R2)<=20/5+15/5*5,5,)))))"
$ws.Range("M2:M6").Formula = '=IF(L2=1,"konservativ",IF(L2=2,"risikoscheu",IF(L2=3,"risikobereit",IF(L2=4,"spekulativ",IF(L2=5,"hochspekulativ")))))'
$ws.Range("L6").Formula = "=IF(SUM(N6:R6)<=20/4+15/5,1,IF(SUM(N6:R6)<=20/4+15/5*2,2,IF(SUM(N6:R6)<=20/4+15/5*3,3,IF(SUM(N6:R6)<=20/4+15/5*4,4,IF(SUM(N6:R6)<=20/4+15/5*5,5,)))))"

# --- 4) add rows 7-11: same five respondents again with ids 5-9 ---
for ($i = 0; $i -lt 5; $i++) {
    $row = 7 + $i
    $id = 5 + $i
    $p = $people[$i]

    $ws.Range("A$row").Value = $id
    $ws.Range("B$row").Value = $p[0]
    $ws.Range("C$row").Value = $p[1]
    $ws.Range("D$row").Value = $p[2]
    $ws.Range("E$row").Value = $p[3]
    $ws.Range("F$row").Value = $p[4]
    $ws.Range("H$row").Value = $p[5]
    $ws.Range("I$row").Value = $p[6]
    $ws.Range("J$row").Value = 35
    $ws.Range("K$row").Value = $p[7]

    $a = $answers[$i]
    $ws.Range("N$row").Value = $a[0]
    $ws.Range("O$row").Value = $a[1]
    $ws.Range("P$row").Value = $a[2]
    $ws.Range("Q$row").Value = $a[3]
    $ws.Range("R$row").Value = $a[4]
}

# same shared-formula layout for the second block of five rows
$ws.Range("L7:L10").Formula = "=IF(SUM(N7:R7)<=20/5+15/5,1,IF(SUM(N7:R7)<=20/5+15/5*2,2,IF(SUM(N7:R7)<=20/5+15/5*3,3,IF(SUM(N7:R7)<=20/5+15/5*4,4,IF(SUM(N7:R7)<=20/5+15/5*5,5,)))))"
$ws.Range("M7:M11").Formula = '=IF(L7=1,"konservativ",IF(L7=2,"risikoscheu",IF(L7=3,"risikobereit",IF(L7=4,"spekulativ",IF(L7=5,"hochspekulativ")))))'
$ws.Range("L11").Formula = "=IF(SUM(N11:R11)<=20/4+15/5,1,IF(SUM(N11:R11)<=20/4+15/5*2,2,IF(SUM(N11:R11)<=20/4+15/5*3,3,IF(SUM(N11:R11)<=20/4+15/5*4,4,IF(SUM(N11:R11)<=20/4+15/5*5,5,)))))"

# final selection sits below the new data, matching the authored sheet
$ws.Range("A12").Select() | Out-Null
